{"js": "// Replace the date header and each \"A\u00d7B=\" table-cell expression with\n// its updated value (document-generation date bump + new practice\n// problems). Each old string is unique in the document, so a plain\n// search/replace per-pair is safe and keeps existing run formatting\n// (font/size) intact because Range.insertText(\"...\", Replace) only\n// swaps the text of the matched range.\nconst replacements = [\n  [\"2025-04-26 Saturday\", \"2025-04-27 Sunday\"],\n  [\"662\u00d73=\", \"335\u00d79=\"],\n  [\"144\u00d78=\", \"122\u00d79=\"],\n  [\"775\u00d72=\", \"339\u00d72=\"],\n  [\"445\u00d74=\", \"953\u00d74=\"],\n  [\"311\u00d78=\", \"889\u00d76=\"],\n  [\"149\u00d77=\", \"865\u00d77=\"],\n  [\"580\u00d72=\", \"470\u00d75=\"],\n  [\"800\u00d75=\", \"736\u00d75=\"],\n  [\"295\u00d72=\", \"649\u00d76=\"],\n  [\"966\u00d76=\", \"128\u00d72=\"],\n  [\"531\u00d77=\", \"112\u00d77=\"],\n  [\"951\u00d78=\", \"214\u00d78=\"],\n  [\"787\u00d73=\", \"239\u00d79=\"],\n  [\"507\u00d72=\", \"984\u00d78=\"],\n  [\"776\u00d72=\", \"746\u00d79=\"],\n  [\"899\u00d75=\", \"358\u00d75=\"],\n  [\"999\u00d73=\", \"821\u00d76=\"],\n  [\"119\u00d79=\", \"862\u00d74=\"],\n  [\"419\u00d72=\", \"309\u00d74=\"],\n  [\"113\u00d79=\", \"490\u00d78=\"],\n  [\"687\u00d79=\", \"390\u00d78=\"],\n  [\"509\u00d75=\", \"140\u00d73=\"],\n  [\"698\u00d78=\", \"539\u00d76=\"],\n  [\"757\u00d72=\", \"984\u00d74=\"],\n  [\"267\u00d73=\", \"338\u00d73=\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const rng of results.items) {\n    rng.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the date header and each \"A\u00d7B=\" table-cell expression with its\n# updated value (document-generation date bump + new practice problems).\n# Each old string occurs exactly once in the whole document, so running\n# Find.Execute(..., Replace:=wdReplaceAll (2)) over the full $d.Content\n# range only ever touches its single match, and leaves the matched run's\n# formatting (font/size) untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-04-26 Saturday\", \"2025-04-27 Sunday\"),\n    @(\"662\u00d73=\", \"335\u00d79=\"),\n    @(\"144\u00d78=\", \"122\u00d79=\"),\n    @(\"775\u00d72=\", \"339\u00d72=\"),\n    @(\"445\u00d74=\", \"953\u00d74=\"),\n    @(\"311\u00d78=\", \"889\u00d76=\"),\n    @(\"149\u00d77=\", \"865\u00d77=\"),\n    @(\"580\u00d72=\", \"470\u00d75=\"),\n    @(\"800\u00d75=\", \"736\u00d75=\"),\n    @(\"295\u00d72=\", \"649\u00d76=\"),\n    @(\"966\u00d76=\", \"128\u00d72=\"),\n    @(\"531\u00d77=\", \"112\u00d77=\"),\n    @(\"951\u00d78=\", \"214\u00d78=\"),\n    @(\"787\u00d73=\", \"239\u00d79=\"),\n    @(\"507\u00d72=\", \"984\u00d78=\"),\n    @(\"776\u00d72=\", \"746\u00d79=\"),\n    @(\"899\u00d75=\", \"358\u00d75=\"),\n    @(\"999\u00d73=\", \"821\u00d76=\"),\n    @(\"119\u00d79=\", \"862\u00d74=\"),\n    @(\"419\u00d72=\", \"309\u00d74=\"),\n    @(\"113\u00d79=\", \"490\u00d78=\"),\n    @(\"687\u00d79=\", \"390\u00d78=\"),\n    @(\"509\u00d75=\", \"140\u00d73=\"),\n    @(\"698\u00d78=\", \"539\u00d76=\"),\n    @(\"757\u00d72=\", \"984\u00d74=\"),\n    @(\"267\u00d73=\", \"338\u00d73=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"No match found for '$oldText'\"\n    }\n}\n"}
